# Apply the "stuff at the bottom of the sheets" commit:
#  1. Fill in the missing pair_kind ("generic") for the 4 practice rows.
#  2. Append a new "stim details" block (rows 27-36) describing the
#     video/audio stimulus counts needed per month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. pair_kind column (J) was blank for the practice rows; it's "generic"
#        for every other pair, so fill it in for rows 2-5 too.
$ws.Range("J2:J5").Value = "generic"

# --- 2. New "stim details" section starting at row 27.
$ws.Range("A27").Value = "stim details"

$headers = @("month", "word_type", "need_audio", "need_image", "word", "count", "find images")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(28, $i + 1).Value = $headers[$i]
}

# Rows 29-32: video stimuli needed for months 6/6/7/7
# Rows 33-36: audio stimuli needed for months 6/6/7/7
$data = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)

$row = 29
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
